$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.433.91"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "'1.782.10"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'229.77"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "'0.5848"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.2743"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "'23.13"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'0.06671"
$ws.Range("E10").Value = "  -4.40%  "
$ws.Range("D11").Value = "'0.07526"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "'1.786.84"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "'4.744"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'0.6055"
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").Value = "'2.023.21"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "'74.74"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "'0.000008607"
$ws.Range("E17").Value = "  -11.04%  "
$ws.Range("D18").Value = "'28.397.28"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "'5.363"
$ws.Range("E19").Value = "  -6.11%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'206.43"
$ws.Range("E21").Value = "  -6.44%  "
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "'6.738"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'151.89"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "'8.079"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").Value = "'0.1244"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "'16.22"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'1.403"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").Value = "'0.06114"
$ws.Range("E30").Value = "  -4.64%  "
$ws.Range("D31").Value = "'1.412"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'3.752"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").Value = "'3.754"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'1.663"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").Value = "'1.041"
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("D36").Value = "'0.6340"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").Value = "'2.501"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "'2.686"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01668"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.138.90"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").Value = "'6.266"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").Value = "'0.8732"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "'1.005"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'100.25"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "'1.933.30"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").Value = "'59.54"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").Value = "'8.354"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").Value = "'1.564"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").Value = "'0.05415"
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").Value = "'0.4462"
$ws.Range("E51").Value = "  -2.00%  "
